$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts existing columns C.. to D..
# and automatically adjusts merged cells and the used range/dimension.
$ws.Range("C1").EntireColumn.Insert()

# The old B1/B2 blank placeholder cells (from the merged header area) are no
# longer needed now that the new blank spacer column lives at C1/C2 instead.
$ws.Range("B1").Clear()
$ws.Range("B2").Clear()

# New header for the inserted column
$ws.Range("C3").Value = "Accuracy after attack"

# New "Accuracy after attack" values for rows 4-28
$values = @{
    4  = 94.2
    5  = 5.2
    6  = 0.1
    7  = 0
    8  = 97.09999999999999
    9  = 82.3
    10 = 3
    11 = 0
    12 = 0
    13 = 79
    14 = 6.4
    15 = 0
    16 = 0
    17 = 0
    18 = 77.5
    19 = 91
    20 = 0
    21 = 0.2
    22 = 0.8
    23 = 6.9
    24 = 0.3
    25 = 94.09999999999999
    26 = 6.2
    27 = 0.1
    28 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
